$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 47577.216926574707
$ws.Range("C2").Value = 88404.6640625
$ws.Range("D2").Value = 53.817543029785156

# Row 3
$ws.Range("B3").Value = 46732.700668334961
$ws.Range("C3").Value = 88404.6640625
$ws.Range("D3").Value = 52.862255096435547

# Row 4
$ws.Range("B4").Value = 75536.823303222656
$ws.Range("C4").Value = 88404.6640625
$ws.Range("D4").Value = 85.444389343261719

# Row 5
$ws.Range("B5").Value = 88404.666595458984
$ws.Range("C5").Value = 88404.6640625

# Row 6
$ws.Range("B6").Value = 91665.790893554688
$ws.Range("C6").Value = 88404.6640625
$ws.Range("D6").Value = 103.68886566162109
